$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "P = Tyrel" legend cell (written first so it becomes the first newly-added shared string)
$ws.Range("J9").Value = "P = Tyrel"

# "P" cells (Tyrel-only availability)
$ws.Range("D4").Value = "P"
$ws.Range("D5").Value = "P"
$ws.Range("D6").Value = "P"
$ws.Range("D7").Value = "P"
$ws.Range("B8").Value = "P"
$ws.Range("D8").Value = "P"
$ws.Range("F8").Value = "P"
$ws.Range("B9").Value = "P"
$ws.Range("D9").Value = "P"
$ws.Range("F9").Value = "P"
$ws.Range("B10").Value = "P"
$ws.Range("C10").Value = "P"
$ws.Range("D10").Value = "P"
$ws.Range("F10").Value = "P"
$ws.Range("B11").Value = "P"
$ws.Range("C11").Value = "P"
$ws.Range("D11").Value = "P"
$ws.Range("F11").Value = "P"
$ws.Range("C12").Value = "P"
$ws.Range("C13").Value = "P"
$ws.Range("F16").Value = "P"
$ws.Range("B17").Value = "P"
$ws.Range("D17").Value = "P"
$ws.Range("E17").Value = "P"
$ws.Range("F17").Value = "P"
$ws.Range("B18").Value = "P"
$ws.Range("D18").Value = "P"
$ws.Range("B19").Value = "P"
$ws.Range("D19").Value = "P"
$ws.Range("B20").Value = "P"
$ws.Range("D20").Value = "P"
$ws.Range("B21").Value = "P"
$ws.Range("D21").Value = "P"

# "T,P" cells (both Tristan and Tyrel available) - overwrite former "T" cells / add new ones
$ws.Range("C7").Value = "T,P"
$ws.Range("E7").Value = "T,P"
$ws.Range("C8").Value = "T,P"
$ws.Range("E8").Value = "T,P"
$ws.Range("C9").Value = "T,P"
$ws.Range("E9").Value = "T,P"
$ws.Range("C14").Value = "T,P"
$ws.Range("E14").Value = "T,P"
$ws.Range("C15").Value = "T,P"
$ws.Range("E15").Value = "T,P"
$ws.Range("B16").Value = "T,P"
$ws.Range("D16").Value = "T,P"
$ws.Range("E16").Value = "T,P"

# Update selection to match the saved view
[void]$ws.Range("F18").Select()
